$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: new day-of-work entry (4th day) ---
# Copy formatting from row 5 (odd-row style) into the relevant row 6 cells
# before writing the new values, so the resulting style indices match what
# Excel itself would produce when the user fills the next row down.
$ws.Range("C5").Copy()
$ws.Range("C6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E5").Copy()
$ws.Range("E6").PasteSpecial(-4122)   # xlPasteFormats

# E6 first gets a copy of E5's current text (so both cells briefly share the
# same string), then E6 is edited to the "20" variant and E5 is edited to
# the "24" variant - this mirrors how the description for day 3 (E5) was
# refined to "24" while the brand new day 4 entry (E6) talks about "20".
$ws.Range("E6").Value = $ws.Range("E5").Value2
$ws.Range("E6").Value = "Documentacion de 20 Procediminetos Almacenados de la base de datos ACC MEX"
$ws.Range("E5").Value = "Documentacion de 24 Procediminetos Almacenados de la base de datos ACC MEX"

# Fill in the rest of the new row 6 entry
$ws.Range("C6").Value = 44019
$ws.Range("D6").Value = 6

# --- Update the active selection to reflect where the user ended up ---
[void]$ws.Range("E8").Select()
